$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift data rows up by one: row N (2..10) takes the values that were
# previously in row N+1 (3..11), reflecting the new full ifoCAST series
# evaluation dropping the oldest quarter and adding a new one at the end.
for ($r = 2; $r -le 10; $r++) {
    $src = $r + 1
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($r, $col).Value = $ws.Cells.Item($src, $col).Value2
    }
}

# New values for the last row (11), representing the newest evaluated quarter.
$ws.Range("B11").Value = 0.3270154541542939
$ws.Range("C11").Value = 0.5200289994386147
$ws.Range("D11").Value = 0.3676720554669373
$ws.Range("E11").Value = 0.6063596750006858
$ws.Range("F11").Value = 0.5708906946728048
$ws.Range("G11").Value = 5
